$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 41 and 42 swap content (ranking order changed) with updated price/volume
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.02"
$ws.Range("E41").Value = "  +5.77%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.06"
$ws.Range("E42").Value = "  -2.79%  "

# Price (D) and Volume(1h) (E) refresh for all other rows
# Leading apostrophe forces numeric-looking strings to stay text, matching source formatting.
$ws.Range("D2").Value = "'69.734.75"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'3.804.71"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'614.34"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "'177.26"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").Value = "'3.803.54"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "'6.48"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "'39.82"
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").Value = "'0.0000255"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "'4.442.08"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "'3.805.27"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "'69.809.22"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "'7.56"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("E19").Value = "  -3.84%  "
$ws.Range("D20").Value = "'16.70"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "'509.04"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "'9.60"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "'86.42"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'0.0000145"
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("D27").Value = "'12.73"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").Value = "'10.55"
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "'2.98"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'31.62"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "'6.13"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "'0.142"
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("D39").Value = "'481.06"
$ws.Range("E39").Value = "  +13.40%  "
$ws.Range("D40").Value = "'0.339"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D43").Value = "'49.77"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").Value = "'44.19"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("D45").Value = "'8.58"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "'2.941.30"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "'0.0363"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "'27.30"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'139.10"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'2.44"
$ws.Range("E51").Value = "  -3.23%  "
